# 2017-02-13 snapshot - chunk 30
# Update the STEO Fig19 "U.S. Working Natural Gas in Storage" sheet from the
# January 2017 STEO vintage to the February 2017 STEO vintage: refresh the
# report-title / source-note text, and refresh the forecast-period storage
# values (rows 86-111, i.e. the "B" history/forecast column used by the
# chart's "Storage level" series) plus the two reference-line cells that
# drive the chart's vertical "Forecast" marker (A117:A118).
#
# All of C/D/E/F/G on these rows, and every downstream chart data cache, are
# plain formulas/derived caches over column B, so they recompute on their own
# once column B is rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title and source-note text: January 2017 -> February 2017 -------------
$ws.Range("A2").Value   = "Short-Term Energy Outlook, February 2017"
$ws.Range("A112").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Forecast-period storage levels (column B, rows 86-111) ----------------
$ws.Range("B86").Value  = 3986.88
$ws.Range("B87").Value  = 3323.6785713999998
$ws.Range("B88").Value  = 2697.1974286
$ws.Range("B89").Value  = 2085.8119999999999
$ws.Range("B90").Value  = 1871.7639999999999
$ws.Range("B91").Value  = 2030.222
$ws.Range("B92").Value  = 2381.585
$ws.Range("B93").Value  = 2666.97
$ws.Range("B94").Value  = 2882.5639999999999
$ws.Range("B95").Value  = 3105.7739999999999
$ws.Range("B96").Value  = 3436.54
$ws.Range("B97").Value  = 3730.3310000000001
$ws.Range("B98").Value  = 3656.989
$ws.Range("B99").Value  = 3142.8359999999998
$ws.Range("B100").Value = 2403.0140000000001
$ws.Range("B101").Value = 1847.693
$ws.Range("B102").Value = 1681.11
$ws.Range("B103").Value = 1872.183
$ws.Range("B104").Value = 2259.5880000000002
$ws.Range("B105").Value = 2587.096
$ws.Range("B106").Value = 2834.1570000000002
$ws.Range("B107").Value = 3053.8290000000002
$ws.Range("B108").Value = 3394.0880000000002
$ws.Range("B109").Value = 3688.3380000000002
$ws.Range("B110").Value = 3602.8789999999999
$ws.Range("B111").Value = 3078.2350000000001

# --- Chart "Forecast" vertical-line reference points ------------------------
$ws.Range("A117").Value = 49.5
$ws.Range("A118").Value = 49.5
